# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 8 and row 9) above the existing
# row 8 ("Femacal de La Calera" Papaya data), pushing the old rows 8-21 down
# to rows 10-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before the current row 8. Excel copies formatting
# (e.g. the date style on column D) down from the row above automatically.
$ws.Range("A8:A9").EntireRow.Insert()

# New row 8: Primera quality, 2021-08-05 (serial 44413)
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Femacal de La Calera"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44413
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100108
$ws.Range("H8").Value = "Tropicales y subtropicales"
$ws.Range("I8").Value = 100108004
$ws.Range("J8").Value = "Papaya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("R8").Value = "Provincia del Elquí"
$ws.Range("S8").Value = 1500
$ws.Range("T8").Value = 10

# New row 9: Segunda quality, 2021-08-05 (serial 44413)
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Femacal de La Calera"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44413
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100108
$ws.Range("H9").Value = "Tropicales y subtropicales"
$ws.Range("I9").Value = 100108004
$ws.Range("J9").Value = "Papaya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 58
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("R9").Value = "Provincia del Elquí"
$ws.Range("S9").Value = 1300
$ws.Range("T9").Value = 10
